$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format before assigning numeric-looking strings,
# so Excel keeps them as text instead of converting to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "68.188.37"
$ws.Range("E2").Value = "  +0.51%  "
$ws.Range("D3").Value = "3.834.85"
$ws.Range("E3").Value = "  -0.35%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.31%  "
$ws.Range("D5").Value = "600.02"
$ws.Range("E5").Value = "  +0.40%  "
$ws.Range("D6").Value = "171.56"
$ws.Range("E6").Value = "  +3.46%  "
$ws.Range("D7").Value = "3.834.74"
$ws.Range("E7").Value = "  -0.32%  "
$ws.Range("E8").Value = "  +0.14%  "
$ws.Range("D9").Value = "0.526"
$ws.Range("E9").Value = "  +0.45%  "
$ws.Range("D10").Value = "0.167"
$ws.Range("E10").Value = "  +1.70%  "
$ws.Range("D11").Value = "6.51"
$ws.Range("E11").Value = "  +2.88%  "
$ws.Range("D12").Value = "0.461"
$ws.Range("E12").Value = "  +1.20%  "
$ws.Range("D13").Value = "0.0000282"
$ws.Range("E13").Value = "  +14.14%  "
$ws.Range("D14").Value = "36.88"
$ws.Range("E14").Value = "  +0.40%  "
$ws.Range("D15").Value = "4.496.12"
$ws.Range("E15").Value = "  +0.00%  "
$ws.Range("D16").Value = "3.853.18"
$ws.Range("E16").Value = "  -0.12%  "
$ws.Range("D17").Value = "68.435.06"
$ws.Range("E17").Value = "  +0.76%  "
$ws.Range("D18").Value = "18.38"
$ws.Range("E18").Value = "  +1.83%  "
$ws.Range("D19").Value = "7.44"
$ws.Range("E19").Value = "  +1.51%  "
$ws.Range("D20").Value = "0.111"
$ws.Range("E20").Value = "  +0.78%  "
$ws.Range("D21").Value = "10.89"
$ws.Range("E21").Value = "  -0.16%  "
$ws.Range("D22").Value = "468.84"
$ws.Range("E22").Value = "  +1.20%  "
$ws.Range("D23").Value = "0.730"
$ws.Range("E23").Value = "  +0.39%  "
$ws.Range("D24").Value = "0.0000157"
$ws.Range("E24").Value = "  -2.95%  "
$ws.Range("D25").Value = "83.49"
$ws.Range("E25").Value = "  +0.50%  "
$ws.Range("D26").Value = "2.27"
$ws.Range("E26").Value = "  +1.50%  "
$ws.Range("D27").Value = "12.15"
$ws.Range("E27").Value = "  +0.62%  "
$ws.Range("D28").Value = "10.45"
$ws.Range("E28").Value = "  +5.04%  "
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("D30").Value = "2.93"
$ws.Range("E30").Value = "  -0.27%  "
$ws.Range("D31").Value = "3.997.70"
$ws.Range("E31").Value = "  -0.05%  "
$ws.Range("D32").Value = "7.73"
$ws.Range("E32").Value = "  +0.23%  "
$ws.Range("D33").Value = "2.29"
$ws.Range("E33").Value = "  -0.67%  "
$ws.Range("D34").Value = "31.02"
$ws.Range("E34").Value = "  +0.21%  "
$ws.Range("D35").Value = "9.36"
$ws.Range("E35").Value = "  +0.81%  "
$ws.Range("D36").Value = "3.812.68"
$ws.Range("E36").Value = "  -0.34%  "
$ws.Range("D37").Value = "3.89"
$ws.Range("E37").Value = "  +19.62%  "
$ws.Range("D38").Value = "0.105"
$ws.Range("E38").Value = "  +0.95%  "
$ws.Range("D39").Value = "5.95"
$ws.Range("E39").Value = "  +1.36%  "
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").Value = "0.139"
$ws.Range("E40").Value = "  +0.17%  "
$ws.Range("B41").Value = "Mantle"
$ws.Range("C41").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D41").Value = "1.02"
$ws.Range("E41").Value = "  -0.27%  "
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  +0.27%  "
$ws.Range("D43").Value = "0.318"
$ws.Range("E43").Value = "  +2.37%  "
$ws.Range("D44").Value = "1.98"
$ws.Range("E44").Value = "  +0.36%  "
$ws.Range("E45").Value = "  +0.00%  "
$ws.Range("D46").Value = "8.73"
$ws.Range("E46").Value = "  +3.15%  "
$ws.Range("D47").Value = "417.88"
$ws.Range("E47").Value = "  -1.98%  "
$ws.Range("D48").Value = "0.000291"
$ws.Range("E48").Value = "  +6.61%  "
$ws.Range("D49").Value = "46.69"
$ws.Range("E49").Value = "  -0.91%  "
$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").Value = "0.0359"
$ws.Range("E50").Value = "  +1.68%  "
$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").Value = "141.54"
$ws.Range("E51").Value = "  -1.52%  "

# Restore the original (default) style on column D so no extra
# per-cell style attribute is introduced.
$ws.Range("D2:D51").Style = "Normal"

